# Auto-generated Excel COM-interop script
# Applies cell value updates to match the target diff across 8 sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 5034.0835
$ws.Range("I28").Value = 656
$ws.Range("J28").Value = 18168.334
$ws.Range("K28").Value = 656
$ws.Range("L28").Value = 18168.334
$ws.Range("M28").Value = -171
$ws.Range("N28").Value = -19138.334

$ws.Range("H31").Value = 197.14285
$ws.Range("I31").Value = 30
$ws.Range("K31").Value = 90
$ws.Range("M31").Value = 140

$ws.Range("H38").Value = 2331.8333
$ws.Range("I38").Value = 1364.6666
$ws.Range("J38").Value = 2654.2222
$ws.Range("K38").Value = 4093.9998
$ws.Range("L38").Value = 7962.6666
$ws.Range("M38").Value = -3721.9998
$ws.Range("N38").Value = -8706.6666

$ws.Range("H39").Value = 225.66667
$ws.Range("I39").Value = 51.153847
$ws.Range("K39").Value = 153.461541
$ws.Range("M39").Value = 142.538459

$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H53").Value = 405.84616
$ws.Range("I53").Value = 597
$ws.Range("K53").Value = 597
$ws.Range("M53").Value = 40

$ws.Range("H94").Value = 5196.4287
$ws.Range("I94").Value = 5196.4287
$ws.Range("K94").Value = 5196.4287
$ws.Range("M94").Value = -4745.4287

$ws.Range("H107").Value = 158.25
$ws.Range("I107").Value = 158.25
$ws.Range("K107").Value = 158.25
$ws.Range("M107").Value = 1761.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2862647
$ws.Range("I32").Value = 3121.4666
$ws.Range("K32").Value = 3121.4666
$ws.Range("M32").Value = -2834.4666

$ws.Range("H74").Value = 8190.5713
$ws.Range("I74").Value = 8035.091
$ws.Range("K74").Value = 8035.091
$ws.Range("M74").Value = -7161.091

$ws.Range("H77").Value = 8190.5713
$ws.Range("I77").Value = 8035.091
$ws.Range("K77").Value = 40175.455
$ws.Range("M77").Value = -35807.455

$ws.Range("H102").Value = 5568.75
$ws.Range("I102").Value = 1137.5
$ws.Range("K102").Value = 1137.5
$ws.Range("M102").Value = 484.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 17510442
$ws.Range("I7").Value = 21250038
$ws.Range("J7").Value = 10031250
$ws.Range("K7").Value = 21250038
$ws.Range("L7").Value = 10031250
$ws.Range("M7").Value = -21249925
$ws.Range("N7").Value = -10031476

$ws.Range("H86").Value = 2935.077
$ws.Range("I86").Value = 1927.4736
$ws.Range("J86").Value = 5670
$ws.Range("K86").Value = 1927.4736
$ws.Range("L86").Value = 5670
$ws.Range("M86").Value = -804.4736
$ws.Range("N86").Value = -7916

$ws.Range("H89").Value = 2935.077
$ws.Range("I89").Value = 1927.4736
$ws.Range("J89").Value = 5670
$ws.Range("K89").Value = 9637.368
$ws.Range("L89").Value = 28350
$ws.Range("M89").Value = -4021.368
$ws.Range("N89").Value = -39582

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 694
$ws.Range("I2").Value = 491
$ws.Range("K2").Value = 491
$ws.Range("M2").Value = -378

$ws.Range("H31").Value = 5873.525
$ws.Range("J31").Value = 7270.92
$ws.Range("L31").Value = 7270.92
$ws.Range("N31").Value = -7860.92

$ws.Range("H34").Value = 5873.525
$ws.Range("J34").Value = 7270.92
$ws.Range("L34").Value = 7270.92
$ws.Range("N34").Value = -7674.92

$ws.Range("H35").Value = 730.1539
$ws.Range("I35").Value = 778.5
$ws.Range("J35").Value = 150
$ws.Range("K35").Value = 778.5
$ws.Range("L35").Value = 150
$ws.Range("M35").Value = -484.5
$ws.Range("N35").Value = -738

$ws.Range("H99").Value = 3065.6667
$ws.Range("I99").Value = 2574
$ws.Range("J99").Value = 6999
$ws.Range("K99").Value = 2574
$ws.Range("L99").Value = 6999
$ws.Range("M99").Value = -1076
$ws.Range("N99").Value = -9995

$ws.Range("H106").Value = 25080.25
$ws.Range("J106").Value = 25080.25
$ws.Range("L106").Value = 25080.25
$ws.Range("N106").Value = -27604.25

$ws.Range("H107").Value = 969
$ws.Range("I107").Value = 805
$ws.Range("J107").Value = 1133
$ws.Range("K107").Value = 805
$ws.Range("L107").Value = 1133
$ws.Range("M107").Value = 1115
$ws.Range("N107").Value = -4973

$ws.Range("H126").Value = 3065.6667
$ws.Range("I126").Value = 2574
$ws.Range("J126").Value = 6999
$ws.Range("K126").Value = 7722
$ws.Range("L126").Value = 20997
$ws.Range("M126").Value = -5252
$ws.Range("N126").Value = -25937

$ws.Range("H131").Value = 61206.332
$ws.Range("J131").Value = 61206.332
$ws.Range("L131").Value = 61206.332
$ws.Range("N131").Value = -71286.33199999999

$ws.Range("H132").Value = 4329.385
$ws.Range("I132").Value = 3643.111
$ws.Range("K132").Value = 10929.333
$ws.Range("M132").Value = -8399.332999999999

$ws.Range("H141").Value = 87672.22
$ws.Range("J141").Value = 87672.22
$ws.Range("L141").Value = 87672.22
$ws.Range("N141").Value = -98032.22

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 47.42857
$ws.Range("I2").Value = 64.166664
$ws.Range("J2").Value = 34.875
$ws.Range("K2").Value = 384.999984
$ws.Range("L2").Value = 209.25
$ws.Range("M2").Value = -271.999984
$ws.Range("N2").Value = -435.25

$ws.Range("H121").Value = 300.33334
$ws.Range("J121").Value = 999
$ws.Range("L121").Value = 2997
$ws.Range("N121").Value = -5617

$ws.Range("H134").Value = 4810.875
$ws.Range("I134").Value = 4447.25
$ws.Range("J134").Value = 5174.5
$ws.Range("K134").Value = 13341.75
$ws.Range("L134").Value = 15523.5
$ws.Range("M134").Value = -8271.75
$ws.Range("N134").Value = -25663.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 15000
$ws.Range("J53").Value = 15000
$ws.Range("L53").Value = 15000
$ws.Range("N53").Value = -16262

$ws.Range("H70").Value = 4192.8
$ws.Range("I70").Value = 3546.5557
$ws.Range("K70").Value = 3546.5557
$ws.Range("M70").Value = -3276.5557

$ws.Range("H73").Value = 4192.8
$ws.Range("I73").Value = 3546.5557
$ws.Range("K73").Value = 3546.5557
$ws.Range("M73").Value = -2610.5557

$ws.Range("H102").Value = 1091.1428
$ws.Range("I102").Value = 1317.8182
$ws.Range("K102").Value = 1317.8182
$ws.Range("M102").Value = 304.1818000000001

$ws.Range("H113").Value = 7772.8184
$ws.Range("I113").Value = 1833
$ws.Range("K113").Value = 1833
$ws.Range("M113").Value = 337

$ws.Range("I122").Value = 4199
$ws.Range("J122").Value = 7994
$ws.Range("K122").Value = 12597
$ws.Range("L122").Value = 23982
$ws.Range("M122").Value = -10147
$ws.Range("N122").Value = -28882

$ws.Range("H126").Value = 5820.5
$ws.Range("I126").Value = 5531.4443
$ws.Range("K126").Value = 16594.3329
$ws.Range("M126").Value = -14124.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5454.25
$ws.Range("I7").Value = 3127.6
$ws.Range("J7").Value = 9332
$ws.Range("K7").Value = 3127.6
$ws.Range("L7").Value = 9332
$ws.Range("M7").Value = -3015.6
$ws.Range("N7").Value = -9556

$ws.Range("H16").Value = 135.26666
$ws.Range("I16").Value = 146.92308
$ws.Range("J16").Value = 59.5
$ws.Range("K16").Value = 146.92308
$ws.Range("L16").Value = 59.5
$ws.Range("M16").Value = 23.07692
$ws.Range("N16").Value = -399.5

$ws.Range("H55").Value = 1110.6666
$ws.Range("I55").Value = 1642
$ws.Range("J55").Value = 446.5
$ws.Range("K55").Value = 1642
$ws.Range("L55").Value = 446.5
$ws.Range("M55").Value = -1469
$ws.Range("N55").Value = -792.5

$ws.Range("H82").Value = 3193.7856
$ws.Range("J82").Value = 3935.7
$ws.Range("L82").Value = 3935.7
$ws.Range("N82").Value = -4657.7

$ws.Range("H85").Value = 3193.7856
$ws.Range("J85").Value = 3935.7
$ws.Range("L85").Value = 3935.7
$ws.Range("N85").Value = -6431.7

$ws.Range("H100").Value = 7301.8335
$ws.Range("I100").Value = 4603.6665
$ws.Range("K100").Value = 4603.6665
$ws.Range("M100").Value = -4062.6665

$ws.Range("H126").Value = 5454.25
$ws.Range("I126").Value = 3127.6
$ws.Range("J126").Value = 9332
$ws.Range("K126").Value = 9382.799999999999
$ws.Range("L126").Value = 27996
$ws.Range("M126").Value = -6912.799999999999
$ws.Range("N126").Value = -32936

$ws.Range("H132").Value = 4668
$ws.Range("I132").Value = 4402
$ws.Range("J132").Value = 5200
$ws.Range("K132").Value = 13206
$ws.Range("L132").Value = 15600
$ws.Range("M132").Value = -10676
$ws.Range("N132").Value = -20660

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 5781
$ws.Range("I38").Value = 5700
$ws.Range("J38").Value = 5862
$ws.Range("K38").Value = 5700
$ws.Range("L38").Value = 5862
$ws.Range("M38").Value = -5227
$ws.Range("N38").Value = -6808

$ws.Range("H49").Value = 10000000
$ws.Range("J49").Value = 10000000
$ws.Range("L49").Value = 10000000
$ws.Range("N49").Value = -10000460

$ws.Range("H62").Value = 7857
$ws.Range("I62").Value = 4749.75
$ws.Range("K62").Value = 4749.75
$ws.Range("M62").Value = -4125.75

$ws.Range("H65").Value = 7857
$ws.Range("I65").Value = 4749.75
$ws.Range("K65").Value = 23748.75
$ws.Range("M65").Value = -20628.75

$ws.Range("H122").Value = 1438.2222
$ws.Range("I122").Value = 1434.8572
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 4304.571599999999
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -1854.571599999999
$ws.Range("N122").Value = -9250

$ws.Range("H132").Value = 3427
$ws.Range("I132").Value = 3231
$ws.Range("J132").Value = 4995
$ws.Range("K132").Value = 9693
$ws.Range("L132").Value = 14985
$ws.Range("M132").Value = -7163
$ws.Range("N132").Value = -20045
